# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# price table to the latest scraped snapshot, row by row.
#
# Values that look numeric (e.g. "1.000", "0.9996") are prefixed with a
# leading apostrophe so Excel stores them as literal text instead of
# silently re-parsing them into numbers and dropping the significant
# trailing/leading zeros - exactly like typing them into the grid by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.235.05"
$ws.Range("E2").Value = "  -2.25%  "
$ws.Range("D3").Value = "1.725.18"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("E5").Value = "  -4.16%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4629"
$ws.Range("E7").Value = "  +3.49%  "
$ws.Range("D8").Value = "'0.3455"
$ws.Range("E8").Value = "  -2.66%  "
$ws.Range("D9").Value = "'42.73"
$ws.Range("E9").Value = "  +1.66%  "
$ws.Range("D10").Value = "'0.07292"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("D11").Value = "'1.052"
$ws.Range("E11").Value = "  -3.95%  "
$ws.Range("D12").Value = "'0.9996"
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "'19.90"
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("D14").Value = "'5.877"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").Value = "1.719.99"
$ws.Range("E15").Value = "  -2.99%  "
$ws.Range("D16").Value = "'6.905"
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("D17").Value = "'89.80"
$ws.Range("E17").Value = "  -3.32%  "
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'0.06317"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").Value = "'0.9996"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'16.58"
$ws.Range("E21").Value = "  -3.31%  "
$ws.Range("D22").Value = "'5.652"
$ws.Range("E22").Value = "  -2.51%  "
$ws.Range("D23").Value = "27.260.57"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").Value = "'10.90"
$ws.Range("E24").Value = "  -3.47%  "
$ws.Range("D25").Value = "'2.128"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'155.22"
$ws.Range("E26").Value = "  -4.78%  "
$ws.Range("D27").Value = "'19.43"
$ws.Range("E27").Value = "  -3.88%  "
$ws.Range("D28").Value = "1.920.87"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("D29").Value = "'2.153"
$ws.Range("E29").Value = "  -0.79%  "
$ws.Range("D30").Value = "'119.51"
$ws.Range("E30").Value = "  -4.47%  "
$ws.Range("D31").Value = "'1.039"
$ws.Range("E31").Value = "  -5.69%  "
$ws.Range("D32").Value = "'0.09093"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "'3.588"
$ws.Range("E33").Value = "  -1.36%  "
$ws.Range("D34").Value = "'5.369"
$ws.Range("E34").Value = "  -3.68%  "
$ws.Range("D35").Value = "'0.02214"
$ws.Range("E35").Value = "  -3.43%  "
$ws.Range("D36").Value = "'0.05876"
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("D37").Value = "'11.15"
$ws.Range("E37").Value = "  -5.75%  "
$ws.Range("D38").Value = "'0.2002"
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").Value = "'4.730"
$ws.Range("E39").Value = "  -4.81%  "
$ws.Range("D40").Value = "'0.5980"
$ws.Range("E40").Value = "  -5.33%  "
$ws.Range("D41").Value = "'1.402"
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("D42").Value = "'1.134"
$ws.Range("E42").Value = "  -4.08%  "
$ws.Range("D43").Value = "'7.515"
$ws.Range("E43").Value = "  -5.40%  "
$ws.Range("D44").Value = "'12.66"
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("D45").Value = "'3.602"
$ws.Range("E45").Value = "  -3.55%  "
$ws.Range("D46").Value = "'0.5644"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("D47").Value = "'119.50"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("D48").Value = "'1.876"
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D49").Value = "'0.06672"
$ws.Range("E49").Value = "  -3.34%  "
$ws.Range("D50").Value = "'1.089"
$ws.Range("E50").Value = "  -4.26%  "
$ws.Range("D51").Value = "'0.9993"
$ws.Range("E51").Value = "  -0.11%  "
